$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.700661
$ws.Range("H2").Value = 68.10198299999999
$ws.Range("I2").Value = 0.08615268874617349
$ws.Range("J2").Value = 0.08615268874617349
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 3.221821579973
$ws.Range("R2").Value = 28.996394219757
$ws.Range("S2").Value = 0.0002225374589396369
$ws.Range("T2").Value = 0.0002225374589396369

$ws.Range("G3").Value = 22.700661
$ws.Range("H3").Value = 68.10198299999999
$ws.Range("I3").Value = 0.08615268874617349
$ws.Range("J3").Value = 0.08615268874617349
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 1.745809467979
$ws.Range("R3").Value = 15.712285211811
$ws.Range("S3").Value = 0.0001205864425304588
$ws.Range("T3").Value = 0.0001205864425304588

$ws.Range("G4").Value = 22.700661
$ws.Range("H4").Value = 68.10198299999999
$ws.Range("I4").Value = 0.08615268874617349
$ws.Range("J4").Value = 0.08615268874617349
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 59.33932354805099
$ws.Range("R4").Value = 534.053911932459
$ws.Range("S4").Value = 0.004098682049826885
$ws.Range("T4").Value = 0.004098682049826885

$ws.Range("G5").Value = 22.700661
$ws.Range("H5").Value = 68.10198299999999
$ws.Range("I5").Value = 0.08615268874617349
$ws.Range("J5").Value = 0.08615268874617349
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 1182.982347158849
$ws.Range("R5").Value = 10646.84112442964
$ws.Range("S5").Value = 0.08171088279487651
$ws.Range("T5").Value = 0.08171088279487651

$ws.Range("I6").Value = 0.5030288587986086
$ws.Range("J6").Value = 0.5030288587986087
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 18.81159202589056
$ws.Range("R6").Value = 169.304328233015
$ws.Range("S6").Value = 0.001299353109456143
$ws.Range("T6").Value = 0.001299353109456143

$ws.Range("I7").Value = 0.5030288587986086
$ws.Range("J7").Value = 0.5030288587986087
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.0007040808761221031
$ws.Range("T7").Value = 0.0007040808761221032

$ws.Range("I8").Value = 0.5030288587986086
$ws.Range("J8").Value = 0.5030288587986087
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 346.4708140938117
$ws.Range("R8").Value = 3118.237326844305
$ws.Range("S8").Value = 0.02393141043081298
$ws.Range("T8").Value = 0.02393141043081298

$ws.Range("I9").Value = 0.5030288587986086
$ws.Range("J9").Value = 0.5030288587986087
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 6907.204739987244
$ws.Range("R9").Value = 62164.84265988519
$ws.Range("S9").Value = 0.4770940143822174
$ws.Range("T9").Value = 0.4770940143822174

$ws.Range("G10").Value = 41.94534433333333
$ws.Range("H10").Value = 125.836033
$ws.Range("I10").Value = 0.159189381961201
$ws.Range("J10").Value = 0.159189381961201
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 5.953148921634112
$ws.Range("R10").Value = 53.578340294707
$ws.Range("S10").Value = 0.0004111955304864515
$ws.Range("T10").Value = 0.0004111955304864514

$ws.Range("G11").Value = 41.94534433333333
$ws.Range("H11").Value = 125.836033
$ws.Range("I11").Value = 0.159189381961201
$ws.Range("J11").Value = 0.159189381961201
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 3.225834669517889
$ws.Range("R11").Value = 29.032512025661
$ws.Range("S11").Value = 0.0002228146508100273
$ws.Range("T11").Value = 0.0002228146508100273

$ws.Range("G12").Value = 41.94534433333333
$ws.Range("H12").Value = 125.836033
$ws.Range("I12").Value = 0.159189381961201
$ws.Range("J12").Value = 0.159189381961201
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 109.6447525792343
$ws.Range("R12").Value = 986.8027732131089
$ws.Range("S12").Value = 0.007573375501834119
$ws.Range("T12").Value = 0.007573375501834119

$ws.Range("G13").Value = 41.94534433333333
$ws.Range("H13").Value = 125.836033
$ws.Range("I13").Value = 0.159189381961201
$ws.Range("J13").Value = 0.159189381961201
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 2185.865948653777
$ws.Range("R13").Value = 19672.79353788399
$ws.Range("S13").Value = 0.1509819962780704
$ws.Range("T13").Value = 0.1509819962780704

$ws.Range("G14").Value = 66.302588
$ws.Range("H14").Value = 198.907764
$ws.Range("I14").Value = 0.2516290704940168
$ws.Range("J14").Value = 0.2516290704940168
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 9.410083205350668
$ws.Range("R14").Value = 84.69074884815599
$ws.Range("S14").Value = 0.0006499726794141223
$ws.Range("T14").Value = 0.0006499726794141221

$ws.Range("G15").Value = 66.302588
$ws.Range("H15").Value = 198.907764
$ws.Range("I15").Value = 0.2516290704940168
$ws.Range("J15").Value = 0.2516290704940168
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 5.099044731865333
$ws.Range("R15").Value = 45.891402586788
$ws.Range("S15").Value = 0.0003522008992373697
$ws.Range("T15").Value = 0.0003522008992373697

$ws.Range("G16").Value = 66.302588
$ws.Range("H16").Value = 198.907764
$ws.Range("I16").Value = 0.2516290704940168
$ws.Range("J16").Value = 0.2516290704940168
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 173.314368308708
$ws.Range("R16").Value = 1559.829314778372
$ws.Range("S16").Value = 0.01197115922274984
$ws.Range("T16").Value = 0.01197115922274984

$ws.Range("G17").Value = 66.302588
$ws.Range("H17").Value = 198.907764
$ws.Range("I17").Value = 0.2516290704940168
$ws.Range("J17").Value = 0.2516290704940168
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 3455.176533183159
$ws.Range("R17").Value = 31096.58879864843
$ws.Range("S17").Value = 0.2386557376926155
$ws.Range("T17").Value = 0.2386557376926155
